# Generate Report for Handoff
# The file 8cb9f867-0cf0-4c0d-85e3-55c5778cd867.md has finished translation
# and is now ready for handoff: update its status on the Overview sheet and
# on each language sheet, and stamp the new handoff datetime.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 8cb9f867...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: row for 8cb9f867...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-03 14:34:13"

# --- de-de sheet: row for 8cb9f867...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-03 14:34:26"
